# Formatting excel output order sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the last data row (INCY) - the sheet now covers 20 rows of data (rows 2-21)
$ws.Rows("22").Delete()

# 2) Refresh the "Price" column (E) with updated values
$prices = @{
    2  = 204.24
    3  = 116.96
    4  = 16.5
    5  = 53.58
    6  = 266.98
    7  = 110.03
    8  = 83.315
    9  = 103.96
    10 = 464.05
    11 = 50.86
    12 = 2448.66
    13 = 73.12
    14 = 144.875
    15 = 51.965
    16 = 91.12
    17 = 57.77
    18 = 44.27
    19 = 52.805
    20 = 706.8200000000001
    21 = 95.41
}
foreach ($r in $prices.Keys) {
    $ws.Range("E$r").Value = $prices[$r]
}

# 3) Add the new "Order (> 0 => Buy)" column header in H1, matching the styling
#    already used by the other header cells (bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Order (> 0 => Buy)"

# 4) Populate the new Order column values for each data row
$orders = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 0
}
foreach ($r in $orders.Keys) {
    $ws.Range("H$r").Value = $orders[$r]
}

$excel.CutCopyMode = 0
